$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.191.52'
$ws.Range("E2").Value = '  -2.88%  '

# Row 3
$ws.Range("D3").Value = '1.712.95'
$ws.Range("E3").Value = '  -3.33%  '

# Row 4
$ws.Range("E4").Value = '  +0.27%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.35'
$ws.Range("E5").Value = '  -6.16%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.27%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4729'
$ws.Range("E7").Value = '  +5.54%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3443'
$ws.Range("E8").Value = '  -3.18%  '

# Row 9
$ws.Range("E9").Value = '  +0.15%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07290'
$ws.Range("E10").Value = '  -1.98%  '

# Row 11
$ws.Range("E11").Value = '  -5.81%  '

# Row 12
$ws.Range("E12").Value = '  +0.28%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.85'
$ws.Range("E13").Value = '  -5.44%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.856'
$ws.Range("E14").Value = '  -3.12%  '

# Row 15
$ws.Range("D15").Value = '1.710.79'
$ws.Range("E15").Value = '  -3.29%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.846'

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.77'
$ws.Range("E17").Value = '  -5.05%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001040'
$ws.Range("E18").Value = '  -2.30%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06359'
$ws.Range("E19").Value = '  -1.18%  '

# Row 20
$ws.Range("E20").Value = '  +0.19%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.49'
$ws.Range("E21").Value = '  -3.75%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.596'
$ws.Range("E22").Value = '  -3.19%  '

# Row 23
$ws.Range("D23").Value = '27.227.61'
$ws.Range("E23").Value = '  -2.83%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.80'
$ws.Range("E24").Value = '  -4.39%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.096'
$ws.Range("E25").Value = '  -0.84%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.13'
$ws.Range("E26").Value = '  -5.00%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.84'
$ws.Range("E27").Value = '  -2.62%  '

# Row 28
$ws.Range("D28").Value = '1.909.77'
$ws.Range("E28").Value = '  -2.93%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.088'
$ws.Range("E29").Value = '  -3.45%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.05'
$ws.Range("E30").Value = '  -3.68%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.013'
$ws.Range("E31").Value = '  -8.96%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09218'
$ws.Range("E32").Value = '  +0.34%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.589'
$ws.Range("E33").Value = '  -2.50%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.301'
$ws.Range("E34").Value = '  -6.88%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02193'
$ws.Range("E35").Value = '  -4.29%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05894'

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.05'
$ws.Range("E37").Value = '  -7.05%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2003'
$ws.Range("E38").Value = '  -5.09%  '

# Row 39
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.742'
$ws.Range("E39").Value = '  -4.67%  '

# Row 40
$ws.Range("E40").Value = '  +0.31%  '

# Row 41
$ws.Range("B41").Value = 'WEMIXTOKEN'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.411'
$ws.Range("E41").Value = '  +1.12%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5926'
$ws.Range("E42").Value = '  -6.41%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.112'
$ws.Range("E43").Value = '  -6.18%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.471'
$ws.Range("E44").Value = '  -5.34%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.64'
$ws.Range("E45").Value = '  -5.39%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.565'
$ws.Range("E46").Value = '  -4.82%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5618'
$ws.Range("E47").Value = '  -4.65%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '118.24'
$ws.Range("E48").Value = '  -3.58%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.842'
$ws.Range("E49").Value = '  -6.14%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06631'
$ws.Range("E50").Value = '  -3.71%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.087'
$ws.Range("E51").Value = '  -4.63%  '
